$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit inserts two new line rows ("line7", "line8") right after the
# existing "line6" row, pushing the "extr1".."extr8" rows down by two
# rows (from rows 8-15 to rows 10-17), and refreshes several of the
# from_bus / to_bus / in_service values for the extr* rows along the way.
#
# Final layout for rows 8-17 (columns A..E):
#   row  A   B        C   D   E
#    8   6   line7    14  11  TRUE
#    9   7   line8    16   9  FALSE
#   10   8   extr1     5  12  TRUE
#   11   9   extr2     5   9  TRUE
#   12  10   extr3    10  11  FALSE
#   13  11   extr4     7   8  TRUE
#   14  12   extr5     9  11  TRUE
#   15  13   extr6     7  11  FALSE
#   16  14   extr7     5   7  TRUE
#   17  15   extr8     8   5  FALSE

# Make the two brand-new rows (16 & 17) match the look of the existing
# data rows (bold/centered/bordered "A" column cell from style index 1)
# by copying the formatting of the row above them before filling values.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $false },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $true  },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $false },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $false }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
